$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The edit moves the existing "_GoBack" bookmark from right after the word
# "December" (in the "Iteration Duration" line) down into the table cell
# under "Things That Not Went Well" - specifically into what becomes the
# last (empty) paragraph of that cell, once the paragraph containing the
# "-Unclear user requirements..." bullet is deleted outright.
#
# Word keeps bookmark names unique, so re-adding a bookmark under the same
# name ("_GoBack") both removes the old occurrence and creates the new one -
# a single Bookmarks.Add call takes care of both halves of the diff.
# ---------------------------------------------------------------------------

# Locate the paragraph holding the bullet that must be removed.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("-Unclear user requirements")) {
        $target = $p
    }
}

# Remember where the previous paragraph ends - that position is stable
# across the upcoming delete (everything being removed comes after it).
$prev = $target.Previous()
$prevEnd = $prev.Range.End

# Delete the whole paragraph - text and paragraph mark together (Unit 4 =
# wdParagraph) - which merges it away, leaving the paragraph that used to
# follow it (already empty) as the new last paragraph of the cell.
$target.Range.Delete(4, 1)

# That surviving paragraph is now immediately after $prev.
$pLast = $prev.Next()

# A range touching both the tail of $prev and the (empty) $pLast lets
# Bookmarks.Add anchor the new bookmark inside $pLast, matching the diff.
$rng = $d.Range($prevEnd - 1, $pLast.Range.End)
$d.Bookmarks.Add("_GoBack", $rng) | Out-Null
